# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (TB), C (d2S), D (K), E (IP) for rows 2-8.
# Column G (sum) is recomputed as B+C+D+E for each row.
$data = @{
    2 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697 }
    3 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697 }
    4 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697 }
    5 = @{ B = 0.01293466051926884; C = 0.04071648406533734; D = 0.1494219747398047; E = 0.4942365360607697 }
    6 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 3.537761648806719; E = 0.4942365360607697 }
    7 = @{ B = 0.6606524410359556; C = 1.655778082260271;  D = 3.537761648806719; E = 0.4942365360607697 }
    8 = @{ B = 1.455362044514542; C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.B + $vals.C + $vals.D + $vals.E
}
